$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.878.48'
$ws.Range("E2").Value = '  -0.55%  '

$ws.Range("D3").Value = '3.062.75'
$ws.Range("E3").Value = '  -0.84%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = "'536.73"
$ws.Range("E5").Value = '  -3.18%  '

$ws.Range("D6").Value = "'133.02"
$ws.Range("E6").Value = '  -3.00%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").Value = '3.057.59'
$ws.Range("E8").Value = '  -0.64%  '

$ws.Range("E9").Value = '  -0.27%  '

$ws.Range("E10").Value = '  -3.89%  '

$ws.Range("D11").Value = "'6.05"
$ws.Range("E11").Value = '  -8.43%  '

$ws.Range("D12").Value = "'0.451"
$ws.Range("E12").Value = '  -0.54%  '

$ws.Range("D13").Value = "'0.0000223"
$ws.Range("E13").Value = '  +2.68%  '

$ws.Range("D14").Value = "'34.18"
$ws.Range("E14").Value = '  -2.43%  '

$ws.Range("D15").Value = '3.559.68'
$ws.Range("E15").Value = '  -0.58%  '

$ws.Range("D16").Value = '62.890.92'
$ws.Range("E16").Value = '  -0.60%  '

$ws.Range("E17").Value = '  -0.06%  '

$ws.Range("D18").Value = '3.067.83'
$ws.Range("E18").Value = '  -0.79%  '

$ws.Range("D19").Value = "'6.62"
$ws.Range("E19").Value = '  -0.43%  '

$ws.Range("D20").Value = "'481.42"
$ws.Range("E20").Value = '  -3.89%  '

$ws.Range("D21").Value = "'13.30"
$ws.Range("E21").Value = '  -1.69%  '

$ws.Range("D22").Value = "'0.694"
$ws.Range("E22").Value = '  -1.63%  '

$ws.Range("D23").Value = "'7.08"
$ws.Range("E23").Value = '  -2.61%  '

$ws.Range("D24").Value = "'78.94"
$ws.Range("E24").Value = '  +1.17%  '

$ws.Range("D25").Value = "'12.07"
$ws.Range("E25").Value = '  -1.90%  '

$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("D27").Value = "'2.69"
$ws.Range("E27").Value = '  -2.31%  '

$ws.Range("D28").Value = "'8.07"
$ws.Range("E28").Value = '  -1.00%  '

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = '  -0.04%  '

$ws.Range("D30").Value = "'25.93"
$ws.Range("E30").Value = '  -1.09%  '

$ws.Range("E31").Value = '  -7.53%  '

$ws.Range("D32").Value = "'1.10"
$ws.Range("E32").Value = '  -0.72%  '

$ws.Range("D33").Value = "'2.36"
$ws.Range("E33").Value = '  -6.29%  '

$ws.Range("D34").Value = "'57.08"
$ws.Range("E34").Value = '  -3.30%  '

$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = "'6.00"
$ws.Range("E35").Value = '  +2.22%  '

$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = "'5.32"
$ws.Range("E36").Value = '  +3.37%  '

$ws.Range("D37").Value = "'483.57"
$ws.Range("E37").Value = '  -9.15%  '

$ws.Range("D38").Value = '3.133.61'
$ws.Range("E38").Value = '  +2.47%  '

$ws.Range("E39").Value = '  -4.63%  '

$ws.Range("D40").Value = "'0.0793"
$ws.Range("E40").Value = '  +0.15%  '

$ws.Range("D41").Value = "'0.115"
$ws.Range("E41").Value = '  -4.05%  '

$ws.Range("D42").Value = "'8.07"
$ws.Range("E42").Value = '  +0.11%  '

$ws.Range("D43").Value = "'2.59"
$ws.Range("E43").Value = '  -2.05%  '

$ws.Range("D44").Value = "'0.251"
$ws.Range("E44").Value = '  -0.72%  '

$ws.Range("E45").Value = '  +0.06%  '

$ws.Range("D46").Value = "'121.52"
$ws.Range("E46").Value = '  +1.04%  '

$ws.Range("E47").Value = '  +7.31%  '

$ws.Range("D48").Value = "'2.00"
$ws.Range("E48").Value = '  -2.89%  '

$ws.Range("D49").Value = "'24.32"
$ws.Range("E49").Value = '  +2.02%  '

$ws.Range("E50").Value = '  +1.76%  '

$ws.Range("E51").Value = '  -1.83%  '
